$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BF column holds a "Date" label (as text, e.g. "6-24-2011-12") that was
# off by a day relative to the real NBA game date because of how the stats
# were originally scraped/labelled. Correct it to the real calendar date
# "2012-06-24" for every data row (rows 2-31).
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)   # column BF
    # Force text so Excel doesn't reinterpret the "yyyy-mm-dd"-looking
    # string as a date serial value.
    $cell.NumberFormat = "@"
    $cell.Value = "2012-06-24"
}
